$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.715.67"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").Value = "2.975.17"
$ws.Range("E3").Value = "  -1.58%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'540.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.26%  "

$ws.Range("D6").Value = "'135.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "2.970.69"
$ws.Range("E8").Value = "  -1.54%  "

$ws.Range("D9").Value = "'0.484"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.67%  "

$ws.Range("D10").Value = "'6.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.07%  "

$ws.Range("E11").Value = "  -2.56%  "

$ws.Range("D12").Value = "'0.441"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.08%  "

$ws.Range("D13").Value = "'0.0000217"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.71%  "

$ws.Range("D14").Value = "'33.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.11%  "

$ws.Range("D15").Value = "3.445.23"
$ws.Range("E15").Value = "  -1.93%  "

$ws.Range("D16").Value = "61.762.04"
$ws.Range("E16").Value = "  -0.49%  "

$ws.Range("D17").Value = "'0.107"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.21%  "

$ws.Range("D18").Value = "2.977.60"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("D19").Value = "'6.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.78%  "

$ws.Range("D20").Value = "'463.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.52%  "

$ws.Range("D21").Value = "'13.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").Value = "'0.648"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.01%  "

$ws.Range("D23").Value = "'7.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").Value = "'79.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.77%  "

$ws.Range("D25").Value = "'12.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.01%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Value = "'2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.54%  "

$ws.Range("E28").Value = "  -3.01%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D30").Value = "'1.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.30%  "

$ws.Range("D31").Value = "'25.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.48%  "

$ws.Range("E32").Value = "  -3.25%  "

$ws.Range("D33").Value = "'2.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.43%  "

$ws.Range("D34").Value = "'5.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").Value = "'53.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.03%  "

$ws.Range("D36").Value = "'5.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.93%  "

$ws.Range("D37").Value = "'447.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.38%  "

$ws.Range("D38").Value = "'0.0799"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.28%  "

$ws.Range("D39").Value = "'0.0383"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("D40").Value = "2.932.17"
$ws.Range("E40").Value = "  -9.39%  "

$ws.Range("D41").Value = "'0.113"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.47%  "

$ws.Range("D42").Value = "'7.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.33%  "

$ws.Range("E43").Value = "  -0.62%  "

$ws.Range("D44").Value = "'26.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.67%  "

$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").Value = "'0.244"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("E47").Value = "  -0.47%  "

$ws.Range("D48").Value = "'1.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.53%  "

$ws.Range("D49").Value = "'114.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.42%  "

$ws.Range("D50").Value = "0.0₃0484"
$ws.Range("E50").Value = "  -2.27%  "

$ws.Range("E51").Value = "  -2.31%  "
